$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "brand"
$ws.Range("N2").Value = "'TATA"
$ws.Range("N3").Value = "'TATA"
$ws.Range("N4").Value = "'TATA"

$ws.Range("N6").Select()
